$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of A1 (shared string "AmazonHome" -> "Amazon Home")
$ws.Range("A1").Value = "Amazon Home"

# Format A1:A2 with Times New Roman, size 14
$rng = $ws.Range("A1:A2")
$rng.Font.Name = "Times New Roman"
$rng.Font.Size = 14

# Widen column A and adjust row heights
$ws.Columns.Item(1).ColumnWidth = 15.91
$ws.Rows.Item(1).RowHeight = 17.35
$ws.Rows.Item(2).RowHeight = 17.35

# Move the active selection to C3
$ws.Range("C3").Select()
